# Add a new "PREFIX_newsample" sample row to the Samples sheet, and make the
# Samples sheet the active/selected sheet (activeTab=1, tabSelected on Samples,
# selection moved to A18), matching the author's edit.

$wb = $excel.ActiveWorkbook
$wsSamples = $wb.Worksheets.Item("Samples")

# --- Samples sheet: append a new sample row (row 17) ---
# Columns: A=Sample Name, B=Date Collected, C=Researcher Name, D=Tissue,
#          E=Collection Time, F=Animal ID
$wsSamples.Range("A17").Value = "PREFIX_newsample"

# Keep the date column as literal text (matches the existing text entries in
# column B) instead of letting it be auto-converted to a date serial number.
$wsSamples.Range("B17").NumberFormat = "@"
$wsSamples.Range("B17").Value = "2020-11-19"

$wsSamples.Range("C17").Value = "Xianfeng Zhang"
$wsSamples.Range("D17").Value = "BAT"
$wsSamples.Range("E17").Value = 150
$wsSamples.Range("F17").Value = 971

# --- Make "Samples" the active sheet/tab, with A18 selected ---
$wsSamples.Activate() | Out-Null
$wsSamples.Range("A18").Select() | Out-Null

Write-Output "done"
